$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp shown in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 19:42"

# Refresh country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) and update the couple of
# country labels whose rank changed as a result of the data refresh.

# Row 4
$ws.Range("B4").Value = 7259994
$ws.Range("C4").Value = 15810
$ws.Range("D4").Value = 4500685
$ws.Range("E4").Value = 2550551
$ws.Range("G4").Value = 318
$ws.Range("H4").Value = 208758

# Row 5
$ws.Range("B5").Value = 5979197
$ws.Range("C5").Value = 77626
$ws.Range("D5").Value = 4926268
$ws.Range("E5").Value = 958496
$ws.Range("G5").Value = 1023
$ws.Range("H5").Value = 94433

# Row 6
$ws.Range("B6").Value = 4694648
$ws.Range("C6").Value = 2069
$ws.Range("E6").Value = 512916
$ws.Range("G6").Value = 74
$ws.Range("H6").Value = 140783

# Row 14
$ws.Range("B14").Value = 527446
$ws.Range("C14").Value = 14412
$ws.Range("E14").Value = 400855
$ws.Range("G14").Value = 39
$ws.Range("H14").Value = 31700

# Row 21
$ws.Range("B21").Value = 312966
$ws.Range("C21").Value = 1511
$ws.Range("D21").Value = 274514
$ws.Range("E21").Value = 30523
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = 7929

# Row 25
$ws.Range("B25").Value = 284469
$ws.Range("C25").Value = 763
$ws.Range("E25").Value = 25438
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 9531

# Row 30
$ws.Range("A30").Value = "Ecuador"
$ws.Range("B30").Value = 133981
$ws.Range("C30").Value = 1506
$ws.Range("D30").Value = 102852
$ws.Range("E30").Value = 19856
$ws.Range("G30").Value = 37
$ws.Range("H30").Value = 11273

# Row 31
$ws.Range("A31").Value = "Bolivia"
$ws.Range("B31").Value = 133222
$ws.Range("C31").Value = 604
$ws.Range("D31").Value = 92868
$ws.Range("E31").Value = 32554
$ws.Range("G31").Value = 35
$ws.Range("H31").Value = 7800

# Row 34
$ws.Range("B34").Value = 115241
$ws.Range("C34").Value = 2719
$ws.Range("D34").Value = 94150
$ws.Range("E34").Value = 19050
$ws.Range("G34").Value = 43
$ws.Range("H34").Value = 2041

# Row 52
$ws.Range("B52").Value = 72700
$ws.Range("C52").Value = 527
$ws.Range("D52").Value = 30029
$ws.Range("E52").Value = 41506
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 1165

# Row 62
$ws.Range("B62").Value = 50914
$ws.Range("C62").Value = 160
$ws.Range("D62").Value = 35756
$ws.Range("E62").Value = 13447
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 1711

# Row 73
$ws.Range("A73").Value = "Libano"
$ws.Range("B73").Value = 35242
$ws.Range("C73").Value = 1280
$ws.Range("D73").Value = 15434
$ws.Range("E73").Value = 19468
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 340

# Row 74
$ws.Range("A74").Value = "Irlanda"
$ws.Range("B74").Value = 34560
$ws.Range("C74").Value = 245
$ws.Range("D74").Value = 23364
$ws.Range("E74").Value = 9394
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 1802

# Row 93
$ws.Range("B93").Value = 13606
$ws.Range("C93").Value = 14
$ws.Range("E93").Value = 6006

# Row 94
$ws.Range("B94").Value = 13589
$ws.Range("C94").Value = 44
$ws.Range("E94").Value = 2129

# Row 103
$ws.Range("B103").Value = 9863
$ws.Range("C103").Value = 32
$ws.Range("D103").Value = 9500
$ws.Range("E103").Value = 298

# Row 118
$ws.Range("A118").Value = "Cuba"
$ws.Range("B118").Value = 5412
$ws.Range("C118").Value = 62
$ws.Range("D118").Value = 4732
$ws.Range("E118").Value = 560
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 120

# Row 119
$ws.Range("A119").Value = "Republica de Yibuti"
$ws.Range("B119").Value = 5409
$ws.Range("D119").Value = 5340
$ws.Range("E119").Value = 8
$ws.Range("H119").Value = 61

# Row 120
$ws.Range("A120").Value = "Suazilandia"
$ws.Range("B120").Value = 5399
$ws.Range("D120").Value = 4767
$ws.Range("E120").Value = 524
$ws.Range("H120").Value = 108

# Row 149
$ws.Range("B149").Value = 2601
$ws.Range("C149").Value = 40
$ws.Range("D149").Value = 2156
$ws.Range("E149").Value = 435

# Row 160
$ws.Range("B160").Value = 1684
$ws.Range("C160").Value = 13
$ws.Range("E160").Value = 293

# Row 180
$ws.Range("A180").Value = "Eritrea"
$ws.Range("B180").Value = 375
$ws.Range("C180").Value = 6
$ws.Range("D180").Value = 341
$ws.Range("E180").Value = 34

# Row 181
$ws.Range("A181").Value = "Gibraltar"
$ws.Range("B181").Value = 372
$ws.Range("C181").Value = 8
$ws.Range("D181").Value = 335
$ws.Range("E181").Value = 37

# Row 206
$ws.Range("A206").Value = "Timor Oriental"

# Row 207
$ws.Range("A207").Value = "Santa Lucia"
